# Notas de seguimiento día 05-04-2015
#
# - Update wording of two status labels used in the "Seguimiento" sheet
#   (column F): the shared-string text is edited in place so every cell
#   referencing it (F6/F7 -> "Cuaderno de estudio", F9 -> "En revisión
#   de editor") picks up the new wording automatically.
# - The "En revisión de editor" status is now also used on row 10
#   (F10), which previously had no status set.
# - Row 9's height shrinks now that its text is shorter (no longer
#   needs two wrapped lines).
# - Selection marker moves from F7 to F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the shared-string text in place (keeps the same shared-string
# indices, so every referencing cell updates together).
[void]$ws.Cells.Replace("En revisión por parte de editor", "En revisión de editor")
[void]$ws.Cells.Replace("En manuscrito de autor", "Cuaderno de estudio")

# F10 now carries the "En revisión de editor" status too.
$ws.Range("F10").Value = "En revisión de editor"

# Row 9 no longer wraps onto two lines with the shorter text.
$ws.Rows("9").RowHeight = 16.5

# Track the author's last selected cell.
[void]$ws.Range("F9").Select()
